$wb = $excel.ActiveWorkbook

# --- Rename sheets: "Freq" -> "freq", "Diff" -> "diff" (Sheet3 is untouched) ---
$wb.Worksheets.Item("Freq").Name = "freq"
$wb.Worksheets.Item("Diff").Name = "diff"

$wsFreq = $wb.Worksheets.Item("freq")
$wsDiff = $wb.Worksheets.Item("diff")

# --- "freq" sheet: move the selection from L145 to F34 ---
$wsFreq.Activate()
$wsFreq.Range("F34").Select()

# --- "diff" sheet: scroll the window so row 10 becomes the top-left visible
#     cell, and move the selection from C2:C155 down to the single cell F7 ---
$wsDiff.Activate()
$excel.Goto($wsDiff.Range("A10"), $true)
$wsDiff.Range("F7").Select()
